# Update the "Pais" (COVID-19 countries) dashboard sheet to a newer data
# refresh (source data re-pulled at 09:52 instead of 09:22).
#
# The refresh re-sorted several neighbouring countries (by total cases,
# descending) which shows up as adjacent-row swaps in the data, plus a
# handful of standalone numeric updates, plus a timestamp update in the
# title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title timestamp update (row 1) ---------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 09:52"

# --- Rusia / Corea del Sur swap (rows 20-21) ------------------------
# Rusia moves up to row 20 with freshly refreshed totals; Corea del Sur
# moves down to row 21 keeping the numbers it previously had at row 20.
$ws.Range("A20").Value = "Rusia"
$ws.Range("B20").Value = 11917
$ws.Range("C20").Value = 1786
$ws.Range("D20").Value = 795
$ws.Range("E20").Value = 11028
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 94

$ws.Range("A21").Value = "Corea del Sur"
$ws.Range("B21").Value = 10450
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 7117
$ws.Range("E21").Value = 3125
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 208

# --- Standalone numeric refresh, Noruega (row 26) -------------------
$ws.Range("F26").Value = 70

# --- Standalone numeric refresh, Polonia (row 31) -------------------
$ws.Range("D31").Value = 318
$ws.Range("E31").Value = 5083

# --- Standalone numeric refresh, Moldavia (row 62) ------------------
$ws.Range("D62").Value = 56
$ws.Range("E62").Value = 1204

# --- Armenia / Azerbaiyan swap (rows 70-71) -------------------------
# Armenia moves up to row 70 with freshly refreshed totals; Azerbaiyan
# moves down to row 71 keeping the numbers it previously had at row 70.
$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 937
$ws.Range("C70").Value = 16
$ws.Range("D70").Value = 149
$ws.Range("E70").Value = 777
$ws.Range("F70").Value = 30
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 11

$ws.Range("A71").Value = "Azerbaiyan"
$ws.Range("B71").Value = 926
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 101
$ws.Range("E71").Value = 816
$ws.Range("F71").Value = 27
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 9

# --- Standalone numeric refresh, Kuwait (row 72) --------------------
$ws.Range("D72").Value = 123
$ws.Range("E72").Value = 786

# --- Standalone numeric refresh, Kazajistan (row 76) ----------------
$ws.Range("D76").Value = 61
$ws.Range("E76").Value = 732

# --- Montenegro / Vietnam swap (rows 109-110) -----------------------
# Montenegro moves up to row 109 with freshly refreshed totals; Vietnam
# moves down to row 110 keeping the numbers it previously had at row 109.
$ws.Range("A109").Value = "Montenegro"
$ws.Range("B109").Value = 255
$ws.Range("C109").Value = 3
$ws.Range("D109").Value = 4
$ws.Range("E109").Value = 249
$ws.Range("F109").Value = 7
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 2

$ws.Range("A110").Value = "Vietnam"
$ws.Range("B110").Value = 255
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 144
$ws.Range("E110").Value = 111
$ws.Range("F110").Value = 8
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 0
